$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 847.1579
$ws.Range("I41").Value = 1661.5
$ws.Range("J41").Value = 254.90909
$ws.Range("K41").Value = 1661.5
$ws.Range("L41").Value = 254.90909
$ws.Range("M41").Value = -1221.5
$ws.Range("N41").Value = -1134.90909

$ws.Range("H74").Value = 3511.3704
$ws.Range("I74").Value = 3930.3
$ws.Range("J74").Value = 3264.9412
$ws.Range("K74").Value = 3930.3
$ws.Range("L74").Value = 3264.9412
$ws.Range("M74").Value = -2994.3
$ws.Range("N74").Value = -5136.9412

$ws.Range("H77").Value = 3511.3704
$ws.Range("I77").Value = 3930.3
$ws.Range("J77").Value = 3264.9412
$ws.Range("K77").Value = 19651.5
$ws.Range("L77").Value = 16324.706
$ws.Range("M77").Value = -14971.5
$ws.Range("N77").Value = -25684.706

$ws.Range("H112").Value = 1603
$ws.Range("I112").Value = 2033.3334
$ws.Range("J112").Value = 1368.2727
$ws.Range("K112").Value = 6100.0002
$ws.Range("L112").Value = 4104.8181
$ws.Range("M112").Value = -4992.0002
$ws.Range("N112").Value = -6320.8181

$ws.Range("H137").Value = 1950.1702
$ws.Range("I137").Value = 1793.8125
$ws.Range("J137").Value = 2283.7334
$ws.Range("K137").Value = 5381.4375
$ws.Range("L137").Value = 6851.2002
$ws.Range("M137").Value = -2831.4375
$ws.Range("N137").Value = -11951.2002

$ws.Range("H139").Value = 45780
$ws.Range("J139").Value = 45780
$ws.Range("L139").Value = 45780
$ws.Range("N139").Value = -56060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8077.299
$ws.Range("I32").Value = 2741.3225
$ws.Range("J32").Value = 21310.52
$ws.Range("K32").Value = 2741.3225
$ws.Range("L32").Value = 21310.52
$ws.Range("M32").Value = -2454.3225
$ws.Range("N32").Value = -21884.52

$ws.Range("H61").Value = 1432.6666
$ws.Range("I61").Value = 921.25
$ws.Range("J61").Value = 3770.5715
$ws.Range("K61").Value = 921.25
$ws.Range("L61").Value = 3770.5715
$ws.Range("M61").Value = -709.25
$ws.Range("N61").Value = -4194.5715

$ws.Range("H74").Value = 13407551
$ws.Range("I74").Value = 12163929
$ws.Range("J74").Value = 18520220
$ws.Range("K74").Value = 12163929
$ws.Range("L74").Value = 18520220
$ws.Range("M74").Value = -12163055
$ws.Range("N74").Value = -18521968

$ws.Range("H77").Value = 13407551
$ws.Range("I77").Value = 12163929
$ws.Range("J77").Value = 18520220
$ws.Range("K77").Value = 60819645
$ws.Range("L77").Value = 92601100
$ws.Range("M77").Value = -60815277
$ws.Range("N77").Value = -92609836

$ws.Range("H97").Value = 1946.5652
$ws.Range("I97").Value = 1989.091
$ws.Range("J97").Value = 1011
$ws.Range("K97").Value = 1989.091
$ws.Range("L97").Value = 1011
$ws.Range("M97").Value = -1493.091
$ws.Range("N97").Value = -2003

$ws.Range("H132").Value = 2236.7114
$ws.Range("I132").Value = 1830.6
$ws.Range("J132").Value = 3072.8235
$ws.Range("K132").Value = 5491.799999999999
$ws.Range("L132").Value = 9218.470499999999
$ws.Range("M132").Value = -2961.799999999999
$ws.Range("N132").Value = -14278.4705

$ws.Range("H136").Value = 1432.6666
$ws.Range("I136").Value = 921.25
$ws.Range("J136").Value = 3770.5715
$ws.Range("K136").Value = 2763.75
$ws.Range("L136").Value = 11311.7145
$ws.Range("M136").Value = -213.75
$ws.Range("N136").Value = -16411.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 187.2
$ws.Range("J64").Value = 139.66667
$ws.Range("L64").Value = 139.66667
$ws.Range("N64").Value = -589.6666700000001

$ws.Range("H67").Value = 187.2
$ws.Range("J67").Value = 139.66667
$ws.Range("L67").Value = 139.66667
$ws.Range("N67").Value = -1699.66667

$ws.Range("H94").Value = 17548.385
$ws.Range("I94").Value = 2076.3333
$ws.Range("J94").Value = 52360.5
$ws.Range("K94").Value = 2076.3333
$ws.Range("L94").Value = 52360.5
$ws.Range("M94").Value = -1625.3333
$ws.Range("N94").Value = -53262.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5129981.5
$ws.Range("I31").Value = 5001177
$ws.Range("J31").Value = 5265565.5
$ws.Range("K31").Value = 5001177
$ws.Range("L31").Value = 5265565.5
$ws.Range("M31").Value = -5000882
$ws.Range("N31").Value = -5266155.5

$ws.Range("H34").Value = 5129981.5
$ws.Range("I34").Value = 5001177
$ws.Range("J34").Value = 5265565.5
$ws.Range("K34").Value = 5001177
$ws.Range("L34").Value = 5265565.5
$ws.Range("M34").Value = -5000975
$ws.Range("N34").Value = -5265969.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 144.76
$ws.Range("I38").Value = 135.41667
$ws.Range("J38").Value = 153.38461
$ws.Range("K38").Value = 406.25001
$ws.Range("L38").Value = 460.15383
$ws.Range("M38").Value = -59.25001000000003
$ws.Range("N38").Value = -1154.15383

$ws.Range("H81").Value = 3146
$ws.Range("I81").Value = 300
$ws.Range("J81").Value = 3959.1428
$ws.Range("K81").Value = 900
$ws.Range("L81").Value = 11877.4284
$ws.Range("M81").Value = 223
$ws.Range("N81").Value = -14123.4284

$ws.Range("H84").Value = 3146
$ws.Range("I84").Value = 300
$ws.Range("J84").Value = 3959.1428
$ws.Range("K84").Value = 2700
$ws.Range("L84").Value = 35632.2852
$ws.Range("M84").Value = 2916
$ws.Range("N84").Value = -46864.2852

$ws.Range("H97").Value = 2976944.5
$ws.Range("I97").Value = 4762192.5
$ws.Range("J97").Value = 1531.1111
$ws.Range("K97").Value = 14286577.5
$ws.Range("L97").Value = 4593.3333
$ws.Range("M97").Value = -14286081.5
$ws.Range("N97").Value = -5585.3333

$ws.Range("H98").Value = 409.77777
$ws.Range("I98").Value = 200
$ws.Range("J98").Value = 436
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 1308
$ws.Range("M98").Value = 898
$ws.Range("N98").Value = -4304

$ws.Range("H107").Value = 1037.0465
$ws.Range("I107").Value = 286.66666
$ws.Range("J107").Value = 1158.7297
$ws.Range("K107").Value = 859.9999799999999
$ws.Range("L107").Value = 3476.189100000001
$ws.Range("M107").Value = 1060.00002
$ws.Range("N107").Value = -7316.189100000001

$ws.Range("H137").Value = 6342
$ws.Range("I137").Value = 953.3333
$ws.Range("J137").Value = 8138.222
$ws.Range("K137").Value = 2859.9999
$ws.Range("L137").Value = 24414.666
$ws.Range("M137").Value = 2240.0001
$ws.Range("N137").Value = -34614.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2436.037
$ws.Range("I126").Value = 2570.8333
$ws.Range("J126").Value = 2166.4443
$ws.Range("K126").Value = 7712.499899999999
$ws.Range("L126").Value = 6499.3329
$ws.Range("M126").Value = -5242.499899999999
$ws.Range("N126").Value = -11439.3329

$ws.Range("H132").Value = 2096
$ws.Range("I132").Value = 1502.3914
$ws.Range("J132").Value = 3461.3
$ws.Range("K132").Value = 4507.174199999999
$ws.Range("L132").Value = 10383.9
$ws.Range("M132").Value = -1977.174199999999
$ws.Range("N132").Value = -15443.9

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3555
$ws.Range("I7").Value = 3392.2222
$ws.Range("J7").Value = 3973.5715
$ws.Range("K7").Value = 3392.2222
$ws.Range("L7").Value = 3973.5715
$ws.Range("M7").Value = -3280.2222
$ws.Range("N7").Value = -4197.5715

$ws.Range("H22").Value = 5556355.5
$ws.Range("I22").Value = 16667116
$ws.Range("J22").Value = 975
$ws.Range("K22").Value = 16667116
$ws.Range("L22").Value = 975
$ws.Range("M22").Value = -16666821
$ws.Range("N22").Value = -1565

$ws.Range("H27").Value = 5556355.5
$ws.Range("I27").Value = 16667116
$ws.Range("J27").Value = 975
$ws.Range("K27").Value = 16667116
$ws.Range("L27").Value = 975
$ws.Range("M27").Value = -16667009
$ws.Range("N27").Value = -1189

$ws.Range("H40").Value = 3673.6553
$ws.Range("I40").Value = 3775.5625
$ws.Range("J40").Value = 3548.2307
$ws.Range("K40").Value = 3775.5625
$ws.Range("L40").Value = 3548.2307
$ws.Range("M40").Value = -3639.5625
$ws.Range("N40").Value = -3820.2307

$ws.Range("H126").Value = 3555
$ws.Range("I126").Value = 3392.2222
$ws.Range("J126").Value = 3973.5715
$ws.Range("K126").Value = 10176.6666
$ws.Range("L126").Value = 11920.7145
$ws.Range("M126").Value = -7706.6666
$ws.Range("N126").Value = -16860.7145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2272.1
$ws.Range("I126").Value = 2861.7856
$ws.Range("J126").Value = 896.1667
$ws.Range("K126").Value = 8585.356800000001
$ws.Range("L126").Value = 2688.5001
$ws.Range("M126").Value = -6115.356800000001
$ws.Range("N126").Value = -7628.5001
